$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'89.427.27"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "'3.167.99"
$ws.Range("E3").Value = "  -3.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'212.45"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "'613.84"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").Value = "'0.691"
$ws.Range("E8").Value = "  -4.21%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'3.160.93"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").Value = "'0.575"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E12").Value = "  -5.42%  "
$ws.Range("D13").Value = "'0.0000252"
$ws.Range("E13").Value = "  -6.95%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "'89.433.35"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'3.762.46"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D16").Value = "'5.24"
$ws.Range("E16").Value = "  -5.21%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'32.64"
$ws.Range("E17").Value = "  -6.00%  "
$ws.Range("D18").Value = "'3.174.12"
$ws.Range("E18").Value = "  -3.99%  "
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'13.31"
$ws.Range("E20").Value = "  -6.04%  "
$ws.Range("D21").Value = "'434.91"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("E22").Value = "  +34.83%  "
$ws.Range("D23").Value = "'8.57"
$ws.Range("E23").Value = "  -5.11%  "
$ws.Range("D24").Value = "'5.02"
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("D25").Value = "'5.12"
$ws.Range("E25").Value = "  -4.09%  "
$ws.Range("D26").Value = "'11.58"
$ws.Range("E26").Value = "  -6.84%  "
$ws.Range("D27").Value = "'3.348.59"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").Value = "'74.92"
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'0.165"
$ws.Range("E30").Value = "  -8.15%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'3.98"
$ws.Range("E32").Value = "  +21.43%  "
$ws.Range("D33").Value = "'8.38"
$ws.Range("E33").Value = "  -5.94%  "
$ws.Range("D34").Value = "'531.50"
$ws.Range("E34").Value = "  -7.72%  "
$ws.Range("D35").Value = "'6.97"
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("E36").Value = "  -6.81%  "
$ws.Range("D37").Value = "'1.26"
$ws.Range("E37").Value = "  -9.44%  "
$ws.Range("D38").Value = "'21.90"
$ws.Range("E38").Value = "  -5.55%  "
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.127"
$ws.Range("E41").Value = "  -9.86%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -6.71%  "
$ws.Range("E44").Value = "  -8.75%  "
$ws.Range("D45").Value = "'149.06"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("D46").Value = "'43.58"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("D47").Value = "'171.73"
$ws.Range("E47").Value = "  -5.37%  "
$ws.Range("D48").Value = "'0.123"
$ws.Range("E48").Value = "  -11.11%  "
$ws.Range("D49").Value = "'1.22"
$ws.Range("E49").Value = "  -10.04%  "
$ws.Range("D50").Value = "'4.03"
$ws.Range("E50").Value = "  -5.86%  "
$ws.Range("D51").Value = "'0.605"
$ws.Range("E51").Value = "  -5.33%  "
